$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The account statement previously listed three overdue periods (2504, 2505,
# 2506). The new data only has a single overdue period (2507), so the two
# extra detail rows (17 and 18) are removed entirely and the remaining
# detail row's period label is updated.
$ws.Rows("17:18").Delete()

# Update the remaining period label to the new period.
$ws.Range("E16").Value = "2507"

# Update the total overdue amount (Valor Mora) to reflect the single period.
$ws.Range("E11").Value = 56940

# Update the period count (Cant. Periodos) to reflect the single period.
$ws.Range("F13").Value = 1
